$wb = $excel.ActiveWorkbook

# Add the new "InvalidLogin" worksheet after the existing "ValidLogin" sheet
$validSheet = $wb.Worksheets.Item("ValidLogin")
$newSheet = $wb.Worksheets.Add($null, $validSheet)
$newSheet.Name = "InvalidLogin"

# Populate headers and invalid-login sample data
$newSheet.Range("A1").Value = "username"
$newSheet.Range("B1").Value = "password"
$newSheet.Range("A2").Value = "abcd"
$newSheet.Range("B2").Value = "xyz"

$newSheet.Range("B2").Select()

$newSheet.Activate()
